# Adds the missing "Turb_height" column (H) to the turbine coordinates
# table on the first sheet, filling in the turbine hub-height value for
# every data row (2-89).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header cell (H1): new column title, formatted like the other headers (G1) ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Turb_height"

# --- Data cells (H2:H89): turbine heights, in row order ---
$heights = @(
    25, 25, 12, 12, 15, 15, 15, 15, 15, 15, 15, 15, 15, 15, 15, 15, 15,
    25, 25, 25, 18, 18, 15, 15, 15, 15, 15, 18, 18, 18, 18, 18, 18,
    10, 6, 10, 15, 15, 15, 15, 12, 12, 15, 15, 15, 12, 12, 18, 18,
    9, 9, 19, 19, 12, 12, 15, 15, 15, 18, 18, 18, 18, 9, 9, 9,
    18, 18, 18, 18, 20, 20, 20, 20, 15, 15, 18, 18, 18, 18, 18, 18,
    18, 15, 15, 15, 18, 18, 18
)

for ($i = 0; $i -lt $heights.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $heights[$i]
}

# Restore the selection/zoom to match where it ended up after these edits.
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("H90").Select()
